# Fixed Excel template generator
# The "Modification" sheet described 4 tracked columns (created_by/created_at/
# updated_by/updated_at) with a name row and a type row, but had no row
# describing *what* each column actually records. This inserts that
# descriptive row right under the column-name row (row 3), pushing the
# existing name/type rows and the SQL-builder formula row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 3 - Excel shifts rows 3..5 down to 4..6 and auto-adjusts
# all relative/absolute formula references (e.g. A$3 -> A$4, row5 self-refs ->
# row6) the same way the interactive UI would.
$ws.Rows("3:3").Insert()

# Populate the new descriptive row.
$ws.Range("A3").Value = "User created record"
$ws.Range("B3").Value = "Timestamp of record create"
$ws.Range("C3").Value = "User updated record"
$ws.Range("D3").Value = "Timestamp of record update"
